# Weekly update: insert a new report date (2021-11-05) ahead of the existing
# rows for "Terminal La Palmera de La Serena" / Palta / Hass, pushing all
# later rows down by 3 and appending the displaced block at the sheet's end
# (handled automatically by Excel's row-insert shifting the used range).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 fresh rows right before the current row 412.
$ws.Rows.Item(412).Insert()
$ws.Rows.Item(412).Insert()
$ws.Rows.Item(412).Insert()

# --- Row 412: Hass / Primera -------------------------------------------
$ws.Cells.Item(412, 1).Value  = 8
$ws.Cells.Item(412, 2).Value  = "Terminal La Palmera de La Serena"
$ws.Cells.Item(412, 3).Value  = "Coquimbo"
$ws.Cells.Item(412, 4).Value  = 44505
$ws.Cells.Item(412, 5).Value  = 4
$ws.Cells.Item(412, 6).Value  = "Fruta"
$ws.Cells.Item(412, 7).Value  = 100106
$ws.Cells.Item(412, 8).Value  = "Oleaginosos"
$ws.Cells.Item(412, 9).Value  = 100106002
$ws.Cells.Item(412, 10).Value = "Palta"
$ws.Cells.Item(412, 11).Value = "Hass"
$ws.Cells.Item(412, 12).Value = "Primera"
$ws.Cells.Item(412, 13).Value = 400
$ws.Cells.Item(412, 14).Value = 2300
$ws.Cells.Item(412, 15).Value = 2400
$ws.Cells.Item(412, 16).Value = 2350
$ws.Cells.Item(412, 17).Value = "`$/kilo (en caja de 17 kilos)"
$ws.Cells.Item(412, 18).Value = "Provincia de Limarí"
$ws.Cells.Item(412, 19).Value = 2350
$ws.Cells.Item(412, 20).Value = 1

# --- Row 413: Hass / Segunda --------------------------------------------
$ws.Cells.Item(413, 1).Value  = 8
$ws.Cells.Item(413, 2).Value  = "Terminal La Palmera de La Serena"
$ws.Cells.Item(413, 3).Value  = "Coquimbo"
$ws.Cells.Item(413, 4).Value  = 44505
$ws.Cells.Item(413, 5).Value  = 4
$ws.Cells.Item(413, 6).Value  = "Fruta"
$ws.Cells.Item(413, 7).Value  = 100106
$ws.Cells.Item(413, 8).Value  = "Oleaginosos"
$ws.Cells.Item(413, 9).Value  = 100106002
$ws.Cells.Item(413, 10).Value = "Palta"
$ws.Cells.Item(413, 11).Value = "Hass"
$ws.Cells.Item(413, 12).Value = "Segunda"
$ws.Cells.Item(413, 13).Value = 300
$ws.Cells.Item(413, 14).Value = 2000
$ws.Cells.Item(413, 15).Value = 2100
$ws.Cells.Item(413, 16).Value = 2050
$ws.Cells.Item(413, 17).Value = "`$/kilo (en caja de 17 kilos)"
$ws.Cells.Item(413, 18).Value = "Provincia de Limarí"
$ws.Cells.Item(413, 19).Value = 2050
$ws.Cells.Item(413, 20).Value = 1

# --- Row 414: Hass / Tercera ---------------------------------------------
$ws.Cells.Item(414, 1).Value  = 8
$ws.Cells.Item(414, 2).Value  = "Terminal La Palmera de La Serena"
$ws.Cells.Item(414, 3).Value  = "Coquimbo"
$ws.Cells.Item(414, 4).Value  = 44505
$ws.Cells.Item(414, 5).Value  = 4
$ws.Cells.Item(414, 6).Value  = "Fruta"
$ws.Cells.Item(414, 7).Value  = 100106
$ws.Cells.Item(414, 8).Value  = "Oleaginosos"
$ws.Cells.Item(414, 9).Value  = 100106002
$ws.Cells.Item(414, 10).Value = "Palta"
$ws.Cells.Item(414, 11).Value = "Hass"
$ws.Cells.Item(414, 12).Value = "Tercera"
$ws.Cells.Item(414, 13).Value = 300
$ws.Cells.Item(414, 14).Value = 1600
$ws.Cells.Item(414, 15).Value = 1700
$ws.Cells.Item(414, 16).Value = 1650
$ws.Cells.Item(414, 17).Value = "`$/kilo (en caja de 17 kilos)"
$ws.Cells.Item(414, 18).Value = "Provincia de Limarí"
$ws.Cells.Item(414, 19).Value = 1650
$ws.Cells.Item(414, 20).Value = 1
